$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Main Info")

# FLOW design -> wain IP assignment fix (refactor of flow config function)
$ws.Range("B13").Value = "1.1.1.2"
$ws.Range("B21").Value = "2.2.2.1"
$ws.Range("B7").Value = "SMART"

$ws.Range("D7").Select()
